# cadastro_turmas.xlsx: drop the weekday columns (Segunda..Sexta) from the
# header row and append two new turmas (1BADM, 3AADM) with their codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the weekday header cells C1:G1 (Segunda, Terça, Quarta, Quinta, Sexta)
$ws.Range("C1:G1").ClearContents()

# Append new turma rows
$ws.Range("A5").Value = "1BADM"
$ws.Range("B5").Value = 10293845
$ws.Range("A6").Value = "3AADM"
$ws.Range("B6").Value = 1029485

# Move the active selection to D5, matching the saved view state
$ws.Range("D5").Select()
